$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (B5:AH5): round each measurement to 2 decimal places
# (custom accuracy) - column order B..AH matches dataset columns J1..AH
$ws.Cells.Item(5, 2).Value  = 0.96    # B5
$ws.Cells.Item(5, 3).Value  = 0.38    # C5
$ws.Cells.Item(5, 4).Value  = 0.59    # D5
$ws.Cells.Item(5, 5).Value  = 2.02    # E5
$ws.Cells.Item(5, 6).Value  = 1.31    # F5
$ws.Cells.Item(5, 7).Value  = 0.76    # G5
$ws.Cells.Item(5, 8).Value  = 11.07   # H5
$ws.Cells.Item(5, 9).Value  = 1.16    # I5
$ws.Cells.Item(5, 10).Value = 0.56    # J5
$ws.Cells.Item(5, 11).Value = 0.4     # K5
$ws.Cells.Item(5, 12).Value = 0.82    # L5
$ws.Cells.Item(5, 13).Value = 0.77    # M5
$ws.Cells.Item(5, 14).Value = 0.29    # N5
$ws.Cells.Item(5, 15).Value = 0.75    # O5
$ws.Cells.Item(5, 16).Value = 1.2     # P5
$ws.Cells.Item(5, 17).Value = 0.89    # Q5
$ws.Cells.Item(5, 18).Value = 0.68    # R5
$ws.Cells.Item(5, 19).Value = 0.26    # S5
$ws.Cells.Item(5, 20).Value = 4.31    # T5
$ws.Cells.Item(5, 21).Value = 2.75    # U5
$ws.Cells.Item(5, 22).Value = 0.69    # V5
$ws.Cells.Item(5, 23).Value = 1.85    # W5
$ws.Cells.Item(5, 24).Value = 0.78    # X5
$ws.Cells.Item(5, 25).Value = 0.4     # Y5
$ws.Cells.Item(5, 26).Value = 4.77    # Z5
$ws.Cells.Item(5, 27).Value = 0.61    # AA5
$ws.Cells.Item(5, 28).Value = 0.71    # AB5
$ws.Cells.Item(5, 29).Value = 0.79    # AC5
$ws.Cells.Item(5, 30).Value = 0.78    # AD5
$ws.Cells.Item(5, 31).Value = 0.56    # AE5
$ws.Cells.Item(5, 32).Value = 10.73   # AF5
$ws.Cells.Item(5, 33).Value = 0.27    # AG5
$ws.Cells.Item(5, 34).Value = 0.89    # AH5

# Remove row 6 entirely (data trimmed to 1000 rows) - also shrinks the
# sheet dimension from A1:AH6 down to A1:AH5
$ws.Rows("6:6").Delete()
